# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# handback report has been generated:
#   - The "Status" column (shared by all sheets) changes from
#     "Ready for handoff" to "Handed back: in sync with en-US".
#   - Each language sheet (zh-cn, de-de) gets two new populated columns,
#     "Latest Target File" (F) and "Latest Handback File" (G), each
#     containing a hyperlinked file name, for both data rows.
#   - The "Latest Handback DateTime" column (H) is updated with the new
#     handback timestamps (different per language sheet).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$hyperColor = 15570276  # BGR value for RGB(100,149,237) == #6495ED

# ---------------------------------------------------------------------
# Overview sheet: update the Status column text (B/C columns, rows 2-3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column (C) for both rows
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# New "Latest Target File" (F) / "Latest Handback File" (G) hyperlinked cells
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/a.md", "", "", "a.md")
$wsZh.Range("F2").Font.Underline = $true
$wsZh.Range("F2").Font.Color = $hyperColor

$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c99eb845f0a07f3ceb556803006d18666cdcd04f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZh.Range("G2").Font.Underline = $true
$wsZh.Range("G2").Font.Color = $hyperColor

$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/a.md", "", "", "a.md")
$wsZh.Range("F3").Font.Underline = $true
$wsZh.Range("F3").Font.Color = $hyperColor

$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c99eb845f0a07f3ceb556803006d18666cdcd04f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$wsZh.Range("G3").Font.Underline = $true
$wsZh.Range("G3").Font.Color = $hyperColor

# Latest Handback DateTime (H) updated timestamps
$wsZh.Range("H2").Value = "2016-03-19 08:27:44"
$wsZh.Range("H3").Value = "2016-03-19 08:27:44"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column (C) for both rows
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# New "Latest Target File" (F) / "Latest Handback File" (G) hyperlinked cells
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/a.md", "", "", "a.md")
$wsDe.Range("F2").Font.Underline = $true
$wsDe.Range("F2").Font.Color = $hyperColor

$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/af93be5a9ad7cfc272a310c0f8e7d15b802e5fed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDe.Range("G2").Font.Underline = $true
$wsDe.Range("G2").Font.Color = $hyperColor

$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/97e1b4f2cb1b4b83b16fdcf798046363bf183853/e2e/a.md", "", "", "a.md")
$wsDe.Range("F3").Font.Underline = $true
$wsDe.Range("F3").Font.Color = $hyperColor

$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/af93be5a9ad7cfc272a310c0f8e7d15b802e5fed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$wsDe.Range("G3").Font.Underline = $true
$wsDe.Range("G3").Font.Color = $hyperColor

# Latest Handback DateTime (H) updated timestamps (new/unique value for de-de)
$wsDe.Range("H2").Value = "2016-03-19 08:27:49"
$wsDe.Range("H3").Value = "2016-03-19 08:27:49"
